# server/LISTAS/mi/MANIJA FIJA.xlsx - fix(gui) step 1 and 2
# Bump the sheet's date stamp by one day and update the three price rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: date stamp 45308 -> 45309 (2024-01-17 -> 2024-01-18)
$ws.Range("A1").Value = 45309

# Step 1/2 price updates
$ws.Range("D34").Value = 206.846
$ws.Range("D35").Value = 293.075
$ws.Range("D36").Value = 396.525
